# deaths_demo.xlsx – "Add files via upload" edit
#
# The "ethnicities" sheet already carries a 6th category row (the
# "In-Hospital Deaths" / shared-string index 21 label, despite the name, is
# used here for what is semantically the "Other" ethnicity bucket). The
# "prop" sheet recomputes each ethnicities row as a row-wise percentage
# breakdown (value / row-total * 100) and is missing that same 6th row, so
# we add it and refresh the whole percentage table to stay in sync. Finally
# the active sheet/selection bookkeeping moves from "ethnicities" to "prop".

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("ethnicities")
$ws3 = $wb.Worksheets.Item("prop")

# Recompute B2:G6 on "prop" as each "ethnicities" row's percentage of that
# row's total (B:G), including the new row 6 ("In-Hospital Deaths" label).
for ($r = 2; $r -le 6; $r++) {
    $total = 0
    for ($c = 2; $c -le 7; $c++) {
        $total = $total + $ws2.Cells.Item($r, $c).Value2
    }
    for ($c = 2; $c -le 7; $c++) {
        $ws3.Cells.Item($r, $c).Value = $ws2.Cells.Item($r, $c).Value2 / $total * 100
    }
}

# Carry the row label for the newly-added row across from "ethnicities".
$ws3.Range("A6").Value = $ws2.Range("A6").Value2

# Column A on "prop" widens (to match the "ethnicities" sheet's column A).
$ws3.Columns("A").ColumnWidth = 16.109375

# "ethnicities" is no longer the active tab / its prior cell selection is
# replaced by the full data range; "prop" becomes the active tab with a
# new selected cell.
$ws2.Range("A1:G6").Select()
$ws3.Activate()
$ws3.Range("G16").Select()
